# Update the "tests" results sheet (Лист1 / sheet1.xml) with 2 new bugs found.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Column E ("17 Apr" results) was filled in for rows that didn't have it yet ---
$ws.Range("E3").Value  = "OK"
$ws.Range("E4").Value  = "OK"
$ws.Range("E5").Value  = "OK"
$ws.Range("E6").Value  = "OK"
$ws.Range("E7").Value  = "OK"
$ws.Range("E8").Value  = "OK"
$ws.Range("E11").Value = "OK"
$ws.Range("E12").Value = "OK"
$ws.Range("E14").Value = "OK"

# Row 8 shrank a touch (13.5pt instead of 15pt)
$ws.Rows.Item(8).RowHeight = 13.5

# --- Row 22 (task #21): brand-new bug found ---
$ws.Range("B22").Value = "При удалении аргументов из определения, не удаляются из вызовов => проблемы при переименовании и вообще всех операциях"
$ws.Range("E22").Value = "NG"

# --- Row 10 (task #9): the old bug regressed -- the fix lost a function ---
$ws.Range("B10").Value = "получили из команд некорректный код, не теряем команды при возврате: потеряли функцию"
$ws.Range("C10").Value = "?"
$ws.Range("D10").Value = "?"
$ws.Range("E10").Value = "NG"

# Update the saved cursor/selection position to where the reviewer left off
$ws.Range("E15").Select()
